# Form the consolidated report: fill in the computed "Absent" values
# for the rows where they were missing/blank, and correct the
# already-present-but-wrong ones.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("H12").Value = 1
$ws.Range("H13").Value = 0
